# repull data, push all data, mean calculation
# Update the "dSF" column (F) values to reflect the repulled data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    2  = -1
    9  = -6
    12 = -2
    14 = 8
    16 = 0
    17 = -1
    18 = -7
    21 = 1
    22 = 5
    23 = -1
    26 = -4
    27 = -4
    28 = -6
    29 = 7
    31 = -1
    32 = -2
    34 = -3
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
